$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column B formatting (styles) into column C for all 115 rows
$ws.Range("B1:B115").Copy()
$ws.Range("C1:C115").PasteSpecial(-4122)

# Populate column C values/text to match the new date column
$ws.Cells.Item(1,3).NumberFormat = "@"
$ws.Cells.Item(1,3).Value = "2025/10/27"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Cells.Item(2,3).NumberFormat = "@"
$ws.Cells.Item(2,3).Value = "上证"
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Cells.Item(3,3).Value = 64.56
$ws.Cells.Item(4,3).Value = 3991.35
$ws.Cells.Item(6,3).Value = 50.66
$ws.Cells.Item(7,3).Value = 5650.3
$ws.Cells.Item(9,3).Value = 56.49
$ws.Cells.Item(10,3).Value = 4709.92
$ws.Cells.Item(12,3).Value = 59.32
$ws.Cells.Item(13,3).Value = 7360.02
$ws.Cells.Item(15,3).Value = 29.82
$ws.Cells.Item(16,3).Value = 2748.84
$ws.Cells.Item(18,3).Value = 96.93000000000001
$ws.Cells.Item(19,3).Value = 6791.69
$ws.Cells.Item(21,3).Value = 66.84
$ws.Cells.Item(22,3).Value = 84394.91
$ws.Cells.Item(24,3).Value = 85.68000000000001
$ws.Cells.Item(25,3).Value = 19909.14
$ws.Cells.Item(27,3).Value = 78.48999999999999
$ws.Cells.Item(28,3).Value = 39894.54
$ws.Cells.Item(30,3).Value = 58.08
$ws.Cells.Item(31,3).Value = 5692.72
$ws.Cells.Item(33,3).Value = 10.55
$ws.Cells.Item(34,3).Value = 33846.15
$ws.Cells.Item(36,3).Value = 31.82
$ws.Cells.Item(37,3).Value = 3416.8
$ws.Cells.Item(39,3).Value = 50.56
$ws.Cells.Item(40,3).Value = 3220.52
$ws.Cells.Item(42,3).Value = 19.54
$ws.Cells.Item(43,3).Value = 7479.65
$ws.Cells.Item(45,3).Value = 32.67
$ws.Cells.Item(46,3).Value = 9169.4
$ws.Cells.Item(48,3).Value = 10.3
$ws.Cells.Item(49,3).Value = 13068.87
$ws.Cells.Item(51,3).Value = 24.78
$ws.Cells.Item(52,3).Value = 12608.52
$ws.Cells.Item(54,3).Value = 18.23
$ws.Cells.Item(55,3).Value = 9609.889999999999
$ws.Cells.Item(57,3).Value = 24.89
$ws.Cells.Item(58,3).Value = 16097.1
$ws.Cells.Item(60,3).Value = 32.88
$ws.Cells.Item(61,3).Value = 17526.85
$ws.Cells.Item(63,3).Value = 21.55
$ws.Cells.Item(64,3).Value = 10670.52
$ws.Cells.Item(66,3).Value = 14.64
$ws.Cells.Item(67,3).Value = 9831.940000000001
$ws.Cells.Item(69,3).Value = 21.39
$ws.Cells.Item(70,3).Value = 3286.28
$ws.Cells.Item(72,3).Value = 45.06
$ws.Cells.Item(73,3).Value = 6149.34
$ws.Cells.Item(75,3).Value = 29.08
$ws.Cells.Item(76,3).Value = 9566.48
$ws.Cells.Item(78,3).Value = 18.25
$ws.Cells.Item(79,3).Value = 2473.77
$ws.Cells.Item(81,3).Value = 56.1
$ws.Cells.Item(82,3).Value = 2740.32
$ws.Cells.Item(84,3).Value = 58.79
$ws.Cells.Item(85,3).Value = 2780.47
$ws.Cells.Item(87,3).Value = 52.34
$ws.Cells.Item(88,3).Value = 3897.41
$ws.Cells.Item(90,3).Value = 47.73
$ws.Cells.Item(91,3).Value = 2090.12
$ws.Cells.Item(93,3).Value = 27.69
$ws.Cells.Item(94,3).Value = 13845.15
$ws.Cells.Item(96,3).Value = 88.23
$ws.Cells.Item(97,3).Value = 9721.639999999999
$ws.Cells.Item(99,3).Value = 58.03
$ws.Cells.Item(100,3).Value = 12453
$ws.Cells.Item(102,3).Value = 3.93
$ws.Cells.Item(103,3).Value = 2280.36
$ws.Cells.Item(105,3).Value = 30.94
$ws.Cells.Item(106,3).Value = 898.88
$ws.Cells.Item(108,3).Value = 29.85
$ws.Cells.Item(109,3).Value = 2793.75
$ws.Cells.Item(111,3).Value = 20.63
$ws.Cells.Item(112,3).Value = 4083.23
$ws.Cells.Item(114,3).Value = 29.02
$ws.Cells.Item(115,3).Value = 3463.09
